$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 'Volume 30   Number  50'
$ws.Range("C9").Value = 'Report Covering the Week  12/11/2023  Through  12/17/2023'

# Crime statistics table updates
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E14").Value = -100
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H14").Value = 100
$ws.Range("J14").Value = 8
$ws.Range("K14").Value = 37.5
$ws.Range("C15").NumberFormat = "General"
$ws.Range("C15").Value = "'0"
$ws.Range("M15").Value = -12.5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 700
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 116.666666666667
$ws.Range("I16").Value = 212
$ws.Range("J16").Value = 208
$ws.Range("K16").Value = 1.923076923076
$ws.Range("L16").Value = 21.142857142857
$ws.Range("M16").Value = -25.352112676056
$ws.Range("N16").Value = -76.830601092896
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = -27.906976744186
$ws.Range("I17").Value = 437
$ws.Range("J17").Value = 505
$ws.Range("K17").Value = -13.465346534653
$ws.Range("L17").Value = -14.145383104125
$ws.Range("M17").Value = 45.666666666666
$ws.Range("N17").Value = -53.902953586497
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -31.25
$ws.Range("I18").Value = 166
$ws.Range("J18").Value = 195
$ws.Range("K18").Value = -14.871794871794
$ws.Range("L18").Value = 0.606060606060
$ws.Range("M18").Value = 36.065573770491
$ws.Range("N18").Value = -82.489451476793
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -20
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 33.333333333333
$ws.Range("I19").Value = 428
$ws.Range("J19").Value = 334
$ws.Range("K19").Value = 28.143712574850
$ws.Range("L19").Value = 34.591194968553
$ws.Range("M19").Value = 47.586206896551
$ws.Range("N19").Value = 4.901960784313
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 132
$ws.Range("J20").Value = 72
$ws.Range("K20").Value = 83.333333333333
$ws.Range("L20").Value = 144.444444444444
$ws.Range("M20").Value = 164
$ws.Range("N20").Value = -46.341463414634
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 105
$ws.Range("H21").Value = 7.619047619047
$ws.Range("I21").Value = 1414
$ws.Range("J21").Value = 1347
$ws.Range("K21").Value = 4.974016332590
$ws.Range("L21").Value = 12.669322709163
$ws.Range("M21").Value = 30.202578268876
$ws.Range("N21").Value = -60.744031093836
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Value = "***.*"
$ws.Range("G22").Value = 1
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 20
$ws.Range("F23").Value = 24
$ws.Range("H23").Value = 4.347826086956
$ws.Range("I23").Value = 239
$ws.Range("J23").Value = 231
$ws.Range("K23").Value = 3.463203463203
$ws.Range("L23").Value = 10.648148148148
$ws.Range("M23").Value = 43.975903614457
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = 13.846153846153
$ws.Range("I24").Value = 858
$ws.Range("J24").Value = 798
$ws.Range("K24").Value = 7.518796992481
$ws.Range("L24").Value = 15.013404825737
$ws.Range("M24").Value = 38.164251207729
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 114.285714285714
$ws.Range("F25").Value = 66
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 57.142857142857
$ws.Range("I25").Value = 607
$ws.Range("J25").Value = 581
$ws.Range("K25").Value = 4.475043029259
$ws.Range("L25").Value = 1.675041876046
$ws.Range("M25").Value = -25.612745098039
$ws.Range("C26").NumberFormat = "General"
$ws.Range("C26").Value = "'0"
$ws.Range("F26").Value = 5
$ws.Range("H26").Value = 150
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E27").Value = "***.*"
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 53
$ws.Range("K27").Value = -23.188405797101
$ws.Range("L27").Value = -22.058823529411
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 37
$ws.Range("K28").Value = -21.621621621621
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 31
$ws.Range("K29").Value = -16.129032258064
